$d = $word.ActiveDocument

function Get-ParagraphContaining($text) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($text)) {
            return $p
        }
    }
    return $null
}

function Wrap-BodyXml($bodyInner) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyInner + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData>' +
        '</pkg:part>' +
        '</pkg:package>'
}

function Set-ParagraphXml($paragraph, $innerXml) {
    $paragraph.Range.InsertXML((Wrap-BodyXml $innerXml))
}

# ------------------------------------------------------------------
# 1) Figure caption: "Imagen 1: Estructura de máquina virtual vs Docker"
#    becomes an italic caption split across 4 runs:
#    "Figura" + " 1:" + "Imagen de" + " Estructura de máquina virtual vs Docker"
# ------------------------------------------------------------------
$captionPara = Get-ParagraphContaining "Imagen 1: Estructura de máquina virtual vs Docker"
$captionXml = '<w:p>' +
  '<w:pPr>' +
    '<w:pStyle w:val="normal0"/>' +
    '<w:jc w:val="center"/>' +
    '<w:rPr><w:i/><w:lang w:val="es-AR"/></w:rPr>' +
  '</w:pPr>' +
  '<w:r><w:rPr><w:i/><w:lang w:val="es-AR"/></w:rPr><w:t>Figura</w:t></w:r>' +
  '<w:r><w:rPr><w:i/><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> 1:</w:t></w:r>' +
  '<w:r><w:rPr><w:i/><w:lang w:val="es-AR"/></w:rPr><w:t>Imagen de</w:t></w:r>' +
  '<w:r><w:rPr><w:i/><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> Estructura de máquina virtual vs Docker</w:t></w:r>' +
  '</w:p>'
Set-ParagraphXml $captionPara $captionXml

# ------------------------------------------------------------------
# 2) Merge the proofErr-wrapped "Kernel" run into the surrounding
#    sentence (no visible text change, only XML simplification).
# ------------------------------------------------------------------
$kernelPara = Get-ParagraphContaining "un SO, y que solamente posee su Kernel y bibliotecas esenciales"
$kernelXml = '<w:p>' +
  '<w:pPr><w:pStyle w:val="normal0"/><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:lastRenderedPageBreak/>' +
  '<w:t xml:space="preserve">un SO, y que solamente posee su Kernel y bibliotecas esenciales, y crea un Container para cada proceso, dentro del cual este se ejecuta. Dentro de este Container, además de los procesos, se encuentran las </w:t></w:r>' +
  '<w:r w:rsidR="00BC7B98"><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>dependencias</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> que este necesita para funcionar. </w:t></w:r>' +
  '</w:p>'
Set-ParagraphXml $kernelPara $kernelXml

# ------------------------------------------------------------------
# 3) "Container Image: " heading: merge proofErr-wrapped "Image" run,
#    then add a new descriptive paragraph right after the heading.
# ------------------------------------------------------------------
$containerImagePara = Get-ParagraphContaining "Container Image"
$containerImageXml = '<w:p>' +
  '<w:pPr><w:pStyle w:val="Ttulo2"/><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr>' +
  '<w:bookmarkStart w:id="4" w:name="_sqeevfdcypf3" w:colFirst="0" w:colLast="0"/>' +
  '<w:bookmarkEnd w:id="4"/>' +
  '<w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">Container Image: </w:t></w:r>' +
  '</w:p>'
Set-ParagraphXml $containerImagePara $containerImageXml

$containerImagePara = Get-ParagraphContaining "Container Image"
$containerImagePara.Range.InsertParagraphAfter()
$containerImagePara = Get-ParagraphContaining "Container Image"
$containerImageDescPara = $containerImagePara.Next()
$containerImageDescXml = '<w:p>' +
  '<w:pPr><w:pStyle w:val="normal0"/><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr>' +
  '<w:t xml:space="preserve">En la Container Image se almacena la configuración del Container. Podríamos que el Container es la instancia de la Container Image. En esta configuración se establece el proceso a ejecutar y sus dependencias, debiéndose especificar las versiones de las mismas. Esta configuración se almacena en una estructura de árbol. </w:t>' +
  '</w:r>' +
  '</w:p>'
Set-ParagraphXml $containerImageDescPara $containerImageDescXml

# ------------------------------------------------------------------
# 4) "Docker File:" heading: merge proofErr-wrapped "File" run,
#    then add a new descriptive paragraph right after the heading.
# ------------------------------------------------------------------
$dockerFilePara = Get-ParagraphContaining "Docker File"
$dockerFileXml = '<w:p>' +
  '<w:pPr><w:pStyle w:val="Ttulo2"/><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr>' +
  '<w:bookmarkStart w:id="5" w:name="_d7mp4euq4g8s" w:colFirst="0" w:colLast="0"/>' +
  '<w:bookmarkEnd w:id="5"/>' +
  '<w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Docker File:</w:t></w:r>' +
  '</w:p>'
Set-ParagraphXml $dockerFilePara $dockerFileXml

$dockerFilePara = Get-ParagraphContaining "Docker File"
$dockerFilePara.Range.InsertParagraphAfter()
$dockerFilePara = Get-ParagraphContaining "Docker File"
$dockerFileDescPara = $dockerFilePara.Next()
$dockerFileDescXml = '<w:p>' +
  '<w:pPr><w:pStyle w:val="normal0"/><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr>' +
  '<w:t>Se llama así al archivo que contiene la Container Image. Esta diseñado de forma tal que la Container Image se pueda escribir en un formato fácil de entender.</w:t>' +
  '</w:r>' +
  '</w:p>'
Set-ParagraphXml $dockerFileDescPara $dockerFileDescXml
